$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row at position 8 (pushes existing rows 8-15 down to 9-16)
$ws.Rows.Item(8).Insert()

# Fill in the new row 8 with the I2C pull-up resistor part (R10-R11)
$ws.Range("A8").Value = "R10-R11"
$ws.Range("B8").Value = 2
$ws.Range("C8").Value = "I2C pull-up resistor"
$ws.Range("D8").Value = "0402, 10k, 1%, 1/16W"
$ws.Range("E8").Value = "ROHM Semiconductor"
$ws.Range("F8").Value = "MCR01MRTF1002"
$ws.Range("G8").Value = "http://www.digikey.com/short/3tbm52"

# Extend the autofilter range to cover the new row
$ws.Range("A1:G16").AutoFilter()

# Add a hyperlink on the "Voltage regulator, 5V" row's Link column (now row 11)
$ws.Hyperlinks.Add($ws.Range("G11"), "http://www.digikey.com/short/3t4td1")

# Update selection to match the author's saved cursor position
$ws.Range("A9").Select()
